# Update "想去人数" (want-to-go count) figures in the 江西-漫展信息 workbook
# to reflect the latest generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6516
$ws1.Range("F5").Value = 401
$ws1.Range("F6").Value = 63
$ws1.Range("F10").Value = 85
$ws1.Range("F13").Value = 383
$ws1.Range("F15").Value = 3221
$ws1.Range("F18").Value = 1880

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6516
$ws4.Range("F5").Value = 401
$ws4.Range("F6").Value = 63
$ws4.Range("F8").Value = 2
$ws4.Range("F11").Value = 85
$ws4.Range("F14").Value = 383
$ws4.Range("F16").Value = 3221
$ws4.Range("F19").Value = 1880
